$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 7 new rows (196-202) below the current last row (195), copying
#    row 195's formatting/styles so the new rows inherit the same cell
#    styles (number formats, alignment, etc.) used throughout the table.
# ---------------------------------------------------------------------------
$newRowCount = 7
for ($i = 0; $i -lt $newRowCount; $i++) {
    $ws.Rows("195:195").Copy()
    $ws.Rows("196:196").Insert()
}

# ---------------------------------------------------------------------------
# 2. Fill in the data for the new rows.
# ---------------------------------------------------------------------------
$rows = @(
    @{ Row=196; A=38;   B="Count and Say";                                                        C="#string";                        D="medium"; E=0; F=1; G=20; H=45960; I=45960; Ht=17 },
    @{ Row=197; A=1526; B="Minimum Number of Increments on Subarrays to Form a Target Array";      C="#array #dynamic-programming ";   D="hard";   E=1; F=0; G=50; H=45960; I=45960; Ht=51 },
    @{ Row=198; A=3289; B="The Two Sneaky Numbers of Digitville";                                  C="#array #bit-minipulation #set "; D="easy";   E=1; F=0; G=5;  H=45961; I=45961; Ht=34 },
    @{ Row=199; A=45;   B="Jump Game II";                                                          C="#array #greedy ";                D="medium"; E=0; F=1; G=20; H=45961; I=45961; Ht=17 },
    @{ Row=200; A=1578; B="Minimum Time to Make Rope Colorful";                                    C="#string #array #greedy ";        D="medium"; E=1; F=1; G=20; H=44837; I=45964; Ht=34 },
    @{ Row=201; A=2257; B="Count Unguarded Cells in the Grid";                                     C="#matrix #simulation ";           D="medium"; E=0; F=1; G=20; H=44837; I=45964; Ht=34 },
    @{ Row=202; A=3217; B="Delete Nodes From Linked List Present in Array";                        C="#set #linked-list ";             D="medium"; E=1; F=0; G=8;  H=44837; I=45964; Ht=51 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 1).Value = $r.A
    $ws.Cells.Item($n, 2).Value = $r.B
    $ws.Cells.Item($n, 3).Value = $r.C
    $ws.Cells.Item($n, 4).Value = $r.D
    $ws.Cells.Item($n, 5).Value = $r.E
    $ws.Cells.Item($n, 6).Value = $r.F
    $ws.Cells.Item($n, 7).Value = $r.G
    $ws.Cells.Item($n, 8).Value = $r.H
    $ws.Cells.Item($n, 9).Value = $r.I
    $ws.Rows($n).RowHeight = $r.Ht
}

# ---------------------------------------------------------------------------
# 3. Update the sheet selection / active cell to match the edited workbook.
# ---------------------------------------------------------------------------
$ws.Range("D199").Select()
